# Updates the math-fact table cells (5 columns x 20 rows) in place, cell by
# cell and in document order, so that duplicate problem text (e.g. "89-80=")
# is replaced with the correct distinct target instead of a blanket Find/Replace.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry is the expected old text and its new replacement, in row-major
# order (5 cells per row) matching the table layout.
$replacements = @(
    ,@("62-25=", "22+37=")
    ,@("74+4=", "69+30=")
    ,@("89-80=", "3+48=")
    ,@("53-37=", "62-17=")
    ,@("31+30=", "37-1=")
    ,@("7+51=", "77-43=")
    ,@("51+35=", "96-13=")
    ,@("67-48=", "75-48=")
    ,@("62+18=", "74-20=")
    ,@("35-34=", "92-78=")
    ,@("25+56=", "42-5=")
    ,@("22+18=", "51-24=")
    ,@("19+46=", "37+40=")
    ,@("44-19=", "94-29=")
    ,@("88-36=", "25+0=")
    ,@("85-11=", "21-18=")
    ,@("93-9=", "73-5=")
    ,@("59+5=", "6+66=")
    ,@("51+16=", "94-22=")
    ,@("2+11=", "52+31=")
    ,@("52+41=", "84+0=")
    ,@("73+22=", "44+10=")
    ,@("82-33=", "5+14=")
    ,@("37-25=", "70-56=")
    ,@("3+86=", "47-30=")
    ,@("4+94=", "84-42=")
    ,@("16+78=", "51+7=")
    ,@("72+20=", "46+44=")
    ,@("53+16=", "15+34=")
    ,@("53+24=", "87-45=")
    ,@("42+2=", "31-0=")
    ,@("76-2=", "4+38=")
    ,@("38-26=", "68-22=")
    ,@("62+25=", "88+1=")
    ,@("15+24=", "23+16=")
    ,@("98-27=", "33-22=")
    ,@("1+17=", "22+31=")
    ,@("97-73=", "42+57=")
    ,@("25-17=", "10+45=")
    ,@("79-40=", "79-26=")
    ,@("65+15=", "40-13=")
    ,@("58+20=", "75-22=")
    ,@("46-2=", "70+14=")
    ,@("69-44=", "15+37=")
    ,@("0+75=", "82-16=")
    ,@("36+20=", "81-6=")
    ,@("19-14=", "35+22=")
    ,@("83-37=", "43-39=")
    ,@("89-54=", "80-12=")
    ,@("12+39=", "59-12=")
    ,@("56+34=", "93-56=")
    ,@("5+5=", "54-45=")
    ,@("43-18=", "85-57=")
    ,@("91+6=", "22+25=")
    ,@("51-15=", "86-59=")
    ,@("42+40=", "42+42=")
    ,@("90-22=", "91-11=")
    ,@("90-24=", "43-36=")
    ,@("57+38=", "13+3=")
    ,@("48+23=", "19+80=")
    ,@("83+0=", "13-6=")
    ,@("90-42=", "6+8=")
    ,@("34+49=", "38+20=")
    ,@("58-1=", "63+33=")
    ,@("84-18=", "69-30=")
    ,@("88-60=", "85+13=")
    ,@("50+21=", "21-17=")
    ,@("42+39=", "21+25=")
    ,@("2+76=", "13+25=")
    ,@("22-8=", "7+44=")
    ,@("18+42=", "66+30=")
    ,@("15+60=", "21+70=")
    ,@("68+29=", "4+7=")
    ,@("54-14=", "59+35=")
    ,@("0+29=", "39+55=")
    ,@("79-71=", "64-25=")
    ,@("89-80=", "86-71=")
    ,@("89-60=", "69+1=")
    ,@("46-19=", "34-10=")
    ,@("86-10=", "48-47=")
    ,@("39+23=", "95-57=")
    ,@("2+36=", "41-23=")
    ,@("60+31=", "82-1=")
    ,@("31+57=", "1+45=")
    ,@("64-48=", "71-70=")
    ,@("87-45=", "49-31=")
    ,@("80-0=", "1+32=")
    ,@("67-53=", "67-28=")
    ,@("32-15=", "57-5=")
    ,@("50+39=", "39-22=")
    ,@("25+26=", "23-6=")
    ,@("95-45=", "60-16=")
    ,@("42+43=", "59-18=")
    ,@("15+54=", "89-7=")
    ,@("39-20=", "11-10=")
    ,@("0+93=", "72-59=")
    ,@("80-35=", "95-84=")
    ,@("73-26=", "10+21=")
    ,@("54-46=", "34+0=")
    ,@("26+63=", "0+16=")
)

$cols = 5
$mismatches = 0
for ($i = 0; $i -lt $replacements.Count; $i++) {
    $row = [int]([math]::Floor($i / $cols)) + 1
    $col = ($i % $cols) + 1
    $old = $replacements[$i][0]
    $new = $replacements[$i][1]
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Trim the trailing end-of-cell marker so only the visible text remains.
    [void]$rng.MoveEnd(1, -1)
    if ($rng.Text -ne $old) {
        $mismatches++
        Write-Host "Mismatch at row" $row "col" $col "expected" $old "got" $rng.Text
    }
    $rng.Text = $new
}
Write-Host "Replaced" $replacements.Count "cells with" $mismatches "mismatches."
